# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn/de-de status columns) and on each per-locale
#   sheet's "Status" column.
# - The per-locale "Latest Handback DateTime" gets bumped to the handback run
#   time.
# - The stale "handback file is not latest" Error Detail is cleared now that
#   the handback is in sync.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status cells (E2 / F2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-18 08:51:32"
$zhcn.Range("P2").Value = ""

# de-de detail sheet
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-18 08:51:39"
$dede.Range("P2").Value = ""

# Refresh column widths to fit the newly-generated report content.
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()

$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()

$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
